# The deck's theme (ppt/theme/theme1.xml, currently the "Integral" / "Red
# Violet" design applied to the slide master) is switched over to the
# stock PowerPoint "Office Theme" colour scheme (the scheme that, before
# this edit, only lived in ppt/theme/theme2.xml, used by the notes
# master). We recolour the 12-slot theme colour scheme that backs
# ppt/theme/theme1.xml via the modern ThemeColorScheme object (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink - same order/semantics as
# <a:clrScheme> in the OOXML) so every slide (they all share the single
# slide master/theme) picks up the default Office colours.

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme reaches the one shared theme used by the
# whole deck (single slide master -> theme1.xml), so slide 1 is enough.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB values below are packed the way PowerPoint's ColorFormat.RGB /
# RGBColor.RGB expects: 0xBBGGRR (blue in the high byte), i.e. the
# decimal form of the target sRGB hex triplet read back to front.
$tcs.Item(1).RGB  = 0         # dk1       -> 000000
$tcs.Item(2).RGB  = 16777215  # lt1       -> FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2       -> 44546A
$tcs.Item(4).RGB  = 15132391  # lt2       -> E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1   -> 5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2   -> ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3   -> A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4   -> FFC000
$tcs.Item(9).RGB  = 12874308  # accent5   -> 4472C4
$tcs.Item(10).RGB = 4697456   # accent6   -> 70AD47
$tcs.Item(11).RGB = 12673797  # hlink     -> 0563C1
$tcs.Item(12).RGB = 7491477   # folHlink  -> 954F72
